$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename the header strings: "<name>_old" -> "<name>_FV2310", "<name>_new" -> "<name>_FV2404"
for ($col = 1; $col -le 21; $col++) {
    $cell = $ws.Cells.Item(1, $col)
    $val = $cell.Value()
    if ($val -ne $null) {
        if ($val.EndsWith("_old")) {
            $cell.Value = $val.Substring(0, $val.Length - 4) + "_FV2310"
        } elseif ($val.EndsWith("_new")) {
            $cell.Value = $val.Substring(0, $val.Length - 4) + "_FV2404"
        }
    }
}

# 2. Turn the used range into an Excel table ("Table1"), preserving the header row's
#    existing formatting (bold / grey fill / borders / wrap / centered) instead of
#    letting table-creation capture it into a new dxf.
$headerRange = $ws.Range("A1:U1")
$scratch = $ws.Range("A100:U100")
$headerRange.Copy()
$scratch.PasteSpecial(-4122)

$headerRange.Style = "Normal"

$tableRange = $ws.Range("A1:U61")
$table = $ws.ListObjects.Add(1, $tableRange, $null, 1)
$table.Name = "Table1"

$scratch.Copy()
$headerRange.PasteSpecial(-4122)
$scratch.Clear()

# 3. Freeze the header row (split after row 1).
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
